$d = $word.ActiveDocument

# 1. Skills bullet: drop "shell, " from the list of languages used, e.g.
#    "...python, shell, scala, javascript..." -> "...python, scala, javascript..."
$found1 = $d.Content.Find.Execute(
    "shell, scala",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "scala",
    2
)
if (-not $found1) {
    throw "edit.ps1: could not find the 'shell, scala' text to update"
}

# 2. Remove the whole "大数据主要使用过Spark，HDFS等。" bullet paragraph entirely.
#    Including the paragraph mark (^p) in the search text makes Find consume
#    the trailing end-of-paragraph, so replacing with "" deletes the whole
#    <w:p> rather than leaving an empty bullet behind.
$found2 = $d.Content.Find.Execute(
    "大数据主要使用过Spark，HDFS等。^p",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "",
    2
)
if (-not $found2) {
    throw "edit.ps1: could not find the '大数据...' paragraph to remove"
}
